## dMagMemoryAndProcessTrack.xlsx — add a second worksheet ("Sheet2") with a
## scratch-pad of polynomial-term / sign bookkeeping used while testing the
## maxS equations, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Sheet1 had been left with B11 selected; record the new last-used cell (T8)
# and let it lose the "active tab" flag once Sheet2 is created/activated below.
$sheet1.Range("T8").Select()

# Insert the new sheet immediately after Sheet1.
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Sheet2"

# ---- header labels for the two "term | coeff" mini tables --------------
$ws.Range("D1").Value = 'term'
$ws.Range("E1").Value = 'coeff'
$ws.Range("G1").Value = 'term'
$ws.Range("H1").Value = 'coeff'

# ---- row-of-equation numbering ------------------------------------------
$ws.Range("B2").Value = 'rows of Eqn'
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("M2").Value = 4
$ws.Range("P2").Value = 5
$ws.Range("S2").Value = 6
$ws.Range("V2").Value = 7
$ws.Range("Y2").Value = 8

# ---- term-in-row grid (rows 3-6) ----------------------------------------
$ws.Range("A3").Value = 'term in row'
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = '+'
$ws.Range("D3").Value = 'x^2y^2'
$ws.Range("E3").Value = 'a'
$ws.Range("F3").Value = '+'
$ws.Range("G3").Value = 'x^2'
$ws.Range("H3").Value = 'e'
$ws.Range("I3").Value = '+'
$ws.Range("J3").Value = 'xy^2'
$ws.Range("K3").Value = 'i'
$ws.Range("L3").Value = '+'
$ws.Range("M3").Value = 'x^2y'
$ws.Range("N3").Value = 'm'
$ws.Range("O3").Value = '+'
$ws.Range("P3").Value = 'y'
$ws.Range("Q3").Value = 'q'
$ws.Range("R3").Value = '-'
$ws.Range("S3").Value = 'x^2y^2'
$ws.Range("T3").Value = 'u'
$ws.Range("U3").Value = '-'
$ws.Range("V3").Value = 'x^3y'
$ws.Range("W3").Value = 'ey'
$ws.Range("X3").Value = '+'
$ws.Range("Y3").Value = 'xy'
$ws.Range("Z3").Value = 'B'

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = '+'
$ws.Range("D4").Value = 'x^3y'
$ws.Range("E4").Value = 'b'
$ws.Range("F4").Value = '-'
$ws.Range("G4").Value = 'xy'
$ws.Range("H4").Value = 'f'
$ws.Range("I4").Value = '-'
$ws.Range("J4").Value = 'xy^2'
$ws.Range("K4").Value = 'j'
$ws.Range("L4").Value = '-'
$ws.Range("M4").Value = 'x^3y'
$ws.Range("N4").Value = 'n'
$ws.Range("O4").Value = '-'
$ws.Range("P4").Value = 'y^3'
$ws.Range("Q4").Value = 'r'
$ws.Range("R4").Value = '+'
$ws.Range("S4").Value = 'xy'
$ws.Range("T4").Value = 'v'
$ws.Range("U4").Value = '+'
$ws.Range("V4").Value = 'x^3y'
$ws.Range("W4").Value = 'z'
$ws.Range("X4").Value = '-'
$ws.Range("Y4").Value = 'xy'
$ws.Range("Z4").Value = 'C'

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = '-'
$ws.Range("D5").Value = 'x^3y'
$ws.Range("E5").Value = 'c'
$ws.Range("F5").Value = '+'
$ws.Range("G5").Value = 'xy'
$ws.Range("H5").Value = 'g'
$ws.Range("I5").Value = '-'
$ws.Range("J5").Value = 'x^2y^2'
$ws.Range("K5").Value = 'k'
$ws.Range("L5").Value = '+'
$ws.Range("M5").Value = 'x^3y'
$ws.Range("N5").Value = 'o'
$ws.Range("O5").Value = '+'
$ws.Range("P5").Value = 'y^3'
$ws.Range("Q5").Value = 's'
$ws.Range("R5").Value = '-'
$ws.Range("S5").Value = 'xy^3'
$ws.Range("T5").Value = 'w'
$ws.Range("U5").Value = '-'
$ws.Range("V5").Value = 'x^2'
$ws.Range("W5").Value = 'A'
$ws.Range("X5").Value = '+'
$ws.Range("Y5").Value = 'y^2'
$ws.Range("Z5").Value = 'D'

$ws.Range("B6").Value = 4
$ws.Range("C6").Value = '-'
$ws.Range("D6").Value = 'x^4'
$ws.Range("E6").Value = 'd'
$ws.Range("F6").Value = '-'
$ws.Range("G6").Value = 'y^2'
$ws.Range("H6").Value = 'h'
$ws.Range("I6").Value = '-'
$ws.Range("J6").Value = 'x^2y'
$ws.Range("K6").Value = 'l'
$ws.Range("L6").Value = '+'
$ws.Range("M6").Value = 'x^4'
$ws.Range("N6").Value = 'p'
$ws.Range("O6").Value = '+'
$ws.Range("P6").Value = 'x^2y^2'
$ws.Range("Q6").Value = 't'
$ws.Range("R6").Value = '+'
$ws.Range("S6").Value = 'xy^3'
$ws.Range("T6").Value = 'ex'

# ---- "Grouping" header row of monomials (row 9) -------------------------
$ws.Range("B9").Value = 'Grouping'
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 'x'
$ws.Range("G9").Value = 'xy'
$ws.Range("H9").Value = 'xy^2'
$ws.Range("J9").Value = 'xy^3'
$ws.Range("K9").Value = 'x^2'
$ws.Range("M9").Value = 'x^2y'
$ws.Range("N9").Value = 'x^2y^2'
$ws.Range("P9").Value = 'x^2y^3'
$ws.Range("Q9").Value = 'x^3'
$ws.Range("S9").Value = 'x^3y'
$ws.Range("T9").Value = 'x^3y^2'
$ws.Range("V9").Value = 'x^3y^3'
$ws.Range("W9").Value = 'x^4'
$ws.Range("Y9").Value = 'y'
$ws.Range("Z9").Value = 'y^2'
$ws.Range("AA9").Value = 'y^3'

# ---- "num terms" row (row 10), with a running total in A10 -------------
$ws.Range("A10").Formula = "=SUM(D10:AA10)"
$ws.Range("B10").Value = 'num terms'
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 2
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = 4
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("S10").Value = 6
$ws.Range("T10").Value = 0
$ws.Range("V10").Value = 0
$ws.Range("W10").Value = 2
$ws.Range("Y10").Value = 1
$ws.Range("Z10").Value = 2
$ws.Range("AA10").Value = 2

# ---- "has negative?" row (row 11) ---------------------------------------
$ws.Range("A11").Value = 'has negative?'
$ws.Range("D11").Value = 'NA'
$ws.Range("E11").Value = 'NA'
$ws.Range("G11").Value = '-'
$ws.Range("H11").Value = '-'
$ws.Range("J11").Value = '-'
$ws.Range("K11").Value = '+'
$ws.Range("M11").Value = '-'
$ws.Range("N11").Value = '+'
$ws.Range("P11").Value = 'NA'
$ws.Range("Q11").Value = 'NA'
$ws.Range("S11").Value = '-'
$ws.Range("T11").Value = 'NA'
$ws.Range("V11").Value = 'NA'
$ws.Range("W11").Value = '+'
$ws.Range("Y11").Value = '-'
$ws.Range("Z11").Value = '+'
$ws.Range("AA11").Value = '-'

# ---- quadrant sign table (rows 12-15) -----------------------------------
$ws.Range("B12").Value = '0-90'
$ws.Range("G12").Value = '+'
$ws.Range("H12").Value = '+'
$ws.Range("J12").Value = '+'
$ws.Range("K12").Value = '+'
$ws.Range("M12").Value = '+'
$ws.Range("N12").Value = '+'
$ws.Range("S12").Value = '+'
$ws.Range("W12").Value = '+'
$ws.Range("Y12").Value = '+'
$ws.Range("Z12").Value = '+'
$ws.Range("AA12").Value = '+'

$ws.Range("B13").Value = '90-180'
$ws.Range("G13").Value = '-'
$ws.Range("H13").Value = '-'
$ws.Range("J13").Value = '-'
$ws.Range("K13").Value = '+'
$ws.Range("M13").Value = '+'
$ws.Range("N13").Value = '+'
$ws.Range("S13").Value = '-'
$ws.Range("W13").Value = '+'
$ws.Range("Y13").Value = '+'
$ws.Range("Z13").Value = '+'
$ws.Range("AA13").Value = '+'

$ws.Range("B14").Value = '180-270'
$ws.Range("G14").Value = '+'
$ws.Range("H14").Value = '-'
$ws.Range("J14").Value = '+'
$ws.Range("K14").Value = '+'
$ws.Range("M14").Value = '-'
$ws.Range("N14").Value = '+'
$ws.Range("S14").Value = '+'
$ws.Range("W14").Value = '+'
$ws.Range("Y14").Value = '-'
$ws.Range("Z14").Value = '+'
$ws.Range("AA14").Value = '-'

$ws.Range("B15").Value = '270-360'
$ws.Range("G15").Value = '-'
$ws.Range("H15").Value = '+'
$ws.Range("J15").Value = '-'
$ws.Range("K15").Value = '+'
$ws.Range("M15").Value = '-'
$ws.Range("N15").Value = '+'
$ws.Range("S15").Value = '-'
$ws.Range("W15").Value = '+'
$ws.Range("Y15").Value = '-'
$ws.Range("Z15").Value = '+'
$ws.Range("AA15").Value = '-'

# ---- trailing "+" markers and closing notes -----------------------------
$ws.Range("K16").Value = '+'
$ws.Range("K17").Value = '+'
$ws.Range("K18").Value = '+'

$ws.Range("B21").Value = 'why make the plus minus table above? Because now I will know which terms are positive and which terms are negative.'
$ws.Range("B22").Value = 'If we know all the term constants (see table above), we can know which of the terms are negative and positive so we may set them appropriately equal to one another'

# ---- page setup (match Sheet1's print layout) ---------------------------
$ws.PageSetup.PaperSize = 1
$ws.PageSetup.Orientation = 1
$ws.PageSetup.LeftMargin = 0.7875 * 72
$ws.PageSetup.RightMargin = 0.7875 * 72
$ws.PageSetup.TopMargin = 1.05277777777778 * 72
$ws.PageSetup.BottomMargin = 1.05277777777778 * 72
$ws.PageSetup.HeaderMargin = 0.7875 * 72
$ws.PageSetup.FooterMargin = 0.7875 * 72
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ws.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12Page &P"

# ---- selection + make Sheet2 the active tab -----------------------------
$ws.Range("R3").Select()
$ws.Activate()
